$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "68÷9="
$t.Cell(1,2).Range.Text = "89÷9="
$t.Cell(1,3).Range.Text = "78÷6="
$t.Cell(1,4).Range.Text = "91÷6="
$t.Cell(1,5).Range.Text = "14÷9="
$t.Cell(5,1).Range.Text = "39÷2="
$t.Cell(5,2).Range.Text = "96÷8="
$t.Cell(5,3).Range.Text = "25÷3="
$t.Cell(5,4).Range.Text = "22÷4="
$t.Cell(5,5).Range.Text = "79÷4="
$t.Cell(9,1).Range.Text = "29÷7="
$t.Cell(9,2).Range.Text = "70÷6="
$t.Cell(9,3).Range.Text = "21÷4="
$t.Cell(9,4).Range.Text = "84÷3="
$t.Cell(9,5).Range.Text = "48÷2="
$t.Cell(13,1).Range.Text = "49÷7="
$t.Cell(13,2).Range.Text = "65÷4="
$t.Cell(13,3).Range.Text = "78÷8="
$t.Cell(13,4).Range.Text = "55÷3="
$t.Cell(13,5).Range.Text = "74÷3="
$t.Cell(17,1).Range.Text = "96÷3="
$t.Cell(17,2).Range.Text = "29÷2="
$t.Cell(17,3).Range.Text = "10÷4="
$t.Cell(17,4).Range.Text = "26÷3="
$t.Cell(17,5).Range.Text = "98÷9="
